# Auto-generated: update cryptos Price (D) and Volume(1h) (E) cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.372.13'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -2.00%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.890.22'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.87%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '523.80'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.45%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '141.45'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -7.19%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.549'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -4.63%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.894.22'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -3.10%  '
$ws.Range('E10').Value = '  -6.04%  '
$ws.Range('E11').Value = '  -2.55%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.358'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.53%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.399.54'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.90%  '
$ws.Range('E14').Value = '  +2.18%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '60.408.63'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.05%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '22.66'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -4.62%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.894.19'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.18%  '
$ws.Range('E18').Value = '  -4.84%  '
$ws.Range('E19').Value = '  -4.32%  '
$ws.Range('E20').Value = '  -3.58%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '355.47'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -6.91%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.60'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.41%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.71'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.81%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '64.56'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.77%  '
$ws.Range('E26').Value = '  -4.36%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.179'
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.81'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -5.53%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0836'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -11.79%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').Value = '  -2.89%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.60'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.65%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '150.01'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -6.66%  '
$ws.Range('E35').Value = '  -8.08%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.55'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -6.70%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.993'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -8.11%  '
$ws.Range('E38').Value = '  -6.36%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '37.59'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.29%  '
$ws.Range('E40').Value = '  -5.84%  '
$ws.Range('E41').Value = '  -6.02%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.283.73'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -5.59%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.645'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.83%  '
$ws.Range('E44').Value = '  -2.36%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '20.11'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -9.33%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.90'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.48%  '
$ws.Range('E48').Value = '  -4.88%  '
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0917'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.13%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '246.98'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -8.79%  '
